$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.552.68'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.913.37'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.75%  '
$ws.Range('D5').Value = '325.80'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('D7').Value = '0.4826'
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('D8').Value = '0.4070'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('D9').Value = '0.08144'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = '23.41'
$ws.Range('E11').Value = '  +4.45%  '
$ws.Range('D12').Value = '1.922.73'
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').Value = '7.104'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = '90.44'
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = '0.06781'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').Value = '29.562.26'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').Value = '5.618'
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('D23').Value = '11.84'
$ws.Range('E23').Value = '  +2.79%  '
$ws.Range('D24').Value = '2.175'
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D25').Value = '2.137.41'
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('D26').Value = '154.84'
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('D27').Value = '20.05'
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('D28').Value = '6.306'
$ws.Range('E28').Value = '  +8.25%  '
$ws.Range('E29').Value = '  -1.73%  '
$ws.Range('D30').Value = '119.72'
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('D31').Value = '1.028'
$ws.Range('E31').Value = '  -2.62%  '
$ws.Range('D32').Value = '0.09553'
$ws.Range('D33').Value = '5.517'
$ws.Range('E33').Value = '  +2.41%  '
$ws.Range('D34').Value = '3.562'
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('D35').Value = '1.393'
$ws.Range('E35').Value = '  -2.46%  '
$ws.Range('D36').Value = '0.02268'
$ws.Range('E36').Value = '  +0.43%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').Value = '1.177'
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('D39').Value = '10.80'
$ws.Range('E39').Value = '  +6.33%  '
$ws.Range('D40').Value = '0.5938'
$ws.Range('E40').Value = '  +0.86%  '
$ws.Range('D41').Value = '7.922'
$ws.Range('E41').Value = '  -5.02%  '
$ws.Range('D42').Value = '0.1855'
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('E43').Value = '  -3.31%  '
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('D45').Value = '0.07731'
$ws.Range('E45').Value = '  -3.43%  '
$ws.Range('D46').Value = '12.41'
$ws.Range('E46').Value = '  +1.31%  '
$ws.Range('D47').Value = '0.5567'
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('D49').Value = '115.62'
$ws.Range('E49').Value = '  +1.99%  '
$ws.Range('D50').Value = '72.68'
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('D51').Value = '1.054'
$ws.Range('E51').Value = '  +1.91%  '
